# Auto-generated edit script
# Updates cached market-price / profit figures (columns H-N) across the
# Leve-profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# matching the upstream "chore: update Sheets via scheduled runner" commit.
# Rows are located by sheet + row index + the Leve Item ID in column G,
# to make sure we touch the right record even though row numbers repeat
# across sheets.

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param($Workbook, $SheetName, $Row, $ExpectedItemId, $Values)

    $ws = $Workbook.Worksheets.Item($SheetName)

    # Sanity-check we are editing the row that matches the Leve Item ID
    # (column G) before touching anything.
    $actualItemId = $ws.Cells.Item($Row, 7).Value()
    if ($actualItemId -ne $ExpectedItemId) {
        throw "Row mismatch on $SheetName row $Row : expected G=$ExpectedItemId but found G=$actualItemId"
    }

    $colIndex = @{ H = 8; I = 9; J = 10; K = 11; L = 12; M = 13; N = 14 }

    foreach ($col in $Values.Keys) {
        $cell = $ws.Cells.Item($Row, $colIndex[$col])
        $newValue = $Values[$col]
        if ($null -eq $newValue) {
            # Matches rows where the diff drops the cell entirely rather
            # than writing a 0 / blank value into it.
            $cell.ClearContents()
        } else {
            $cell.Value = $newValue
        }
    }
}


# ALC!6 (Leve Item ID 4564): H6: 38.285713 -> 49.2, I6: 44.333332 -> 49.2, J6: 2 -> 0, K6: 132.999996 -> 147.6, L6: 6 -> 0, M6: -20.99999600000001 -> -35.60000000000002, N6: -230 -> ABSENT
Set-LeveRow $wb "ALC" 6 4564 @{ H = 49.2; I = 49.2; J = 0; K = 147.6; L = 0; M = -35.60000000000002; N = $null }

# ALC!9 (Leve Item ID 5487): H9: 156.66667 -> 175, I9: 147.5 -> 175, K9: 147.5 -> 175, M9: 21.5 -> -6
Set-LeveRow $wb "ALC" 9 5487 @{ H = 175; I = 175; K = 175; M = -6 }

# ALC!11 (Leve Item ID 5533): H11: 275.07144 -> 301.58334, I11: 275.07144 -> 301.58334, K11: 275.07144 -> 301.58334, M11: -135.07144 -> -161.58334
Set-LeveRow $wb "ALC" 11 5533 @{ H = 301.58334; I = 301.58334; K = 301.58334; M = -161.58334 }

# ALC!12 (Leve Item ID 5515): H12: 977.1111 -> 969.3, J12: 933.3333 -> 924.75, L12: 933.3333 -> 924.75, N12: -1273.3333 -> -1264.75
Set-LeveRow $wb "ALC" 12 5515 @{ H = 969.3; J = 924.75; L = 924.75; N = -1264.75 }

# ALC!21 (Leve Item ID 2149): H21: 15000 -> 11499.75, I21: 10000 -> 3999.5, J21: 20000 -> 19000, K21: 10000 -> 3999.5, L21: 20000 -> 19000, M21: -9532 -> -3531.5, N21: -20936 -> -19936
Set-LeveRow $wb "ALC" 21 2149 @{ H = 11499.75; I = 3999.5; J = 19000; K = 3999.5; L = 19000; M = -3531.5; N = -19936 }

# ALC!23 (Leve Item ID 2149): H23: 15000 -> 11499.75, I23: 10000 -> 3999.5, J23: 20000 -> 19000, K23: 10000 -> 3999.5, L23: 20000 -> 19000, M23: -9766 -> -3765.5, N23: -20468 -> -19468
Set-LeveRow $wb "ALC" 23 2149 @{ H = 11499.75; I = 3999.5; J = 19000; K = 3999.5; L = 19000; M = -3765.5; N = -19468 }

# ALC!29 (Leve Item ID 4575): H29: 586.1429000000001 -> 250, I29: 360.6 -> 250, J29: 1150 -> 0, K29: 1081.8 -> 750, L29: 3450 -> 0, M29: -800.8000000000002 -> -469, N29: -4012 -> ABSENT
Set-LeveRow $wb "ALC" 29 4575 @{ H = 250; I = 250; J = 0; K = 750; L = 0; M = -469; N = $null }

# ALC!42 (Leve Item ID 4600): H42: 382.16666 -> 383.83334, I42: 123.5 -> 126, K42: 370.5 -> 378, M42: -140.5 -> -148
Set-LeveRow $wb "ALC" 42 4600 @{ H = 383.83334; I = 126; K = 378; M = -148 }

# ALC!116 (Leve Item ID 27778): H116: 4501.6665 -> 4100, I116: 3500 -> 3498.3333, K116: 3500 -> 3498.3333, M116: -58 -> -56.33329999999978
Set-LeveRow $wb "ALC" 116 27778 @{ H = 4100; I = 3498.3333; K = 3498.3333; M = -56.33329999999978 }

# ALC!132 (Leve Item ID 44049): H132: 48730.57 -> 48728.785, I132: 55984.418 -> 55982.332, K132: 167953.254 -> 167946.996, M132: -165423.254 -> -165416.996
Set-LeveRow $wb "ALC" 132 44049 @{ H = 48728.785; I = 55982.332; K = 167946.996; M = -165416.996 }

# ALC!137 (Leve Item ID 44013): H137: 4333 -> 2824.625, I137: 0 -> 1333, J137: 4333 -> 3719.6, K137: 0 -> 3999, L137: 12999 -> 11158.8, M137: ABSENT -> -1449, N137: -18099 -> -16258.8
Set-LeveRow $wb "ALC" 137 44013 @{ H = 2824.625; I = 1333; J = 3719.6; K = 3999; L = 11158.8; M = -1449; N = -16258.8 }

# ARM!3 (Leve Item ID 2494): H3: 562.5 -> 520.6667, I3: 562.5 -> 520.6667, K3: 562.5 -> 520.6667, M3: -447.5 -> -405.6667
Set-LeveRow $wb "ARM" 3 2494 @{ H = 520.6667; I = 520.6667; K = 520.6667; M = -405.6667 }

# ARM!6 (Leve Item ID 2226): H6: 9999 -> 9997, I6: 9999 -> 9997, K6: 9999 -> 9997, M6: -9826 -> -9824
Set-LeveRow $wb "ARM" 6 2226 @{ H = 9997; I = 9997; K = 9997; M = -9824 }

# ARM!12 (Leve Item ID 2230): H12: 571.4286 -> 600.4, I12: 500 -> 500.5, K12: 500 -> 500.5, M12: -327 -> -327.5
Set-LeveRow $wb "ARM" 12 2230 @{ H = 600.4; I = 500.5; K = 500.5; M = -327.5 }

# ARM!45 (Leve Item ID 27714): H45: 1749 -> 2124, J45: 1500 -> 3000, L45: 1500 -> 3000, N45: -2254 -> -3754
Set-LeveRow $wb "ARM" 45 27714 @{ H = 2124; J = 3000; L = 3000; N = -3754 }

# ARM!61 (Leve Item ID 43999): H61: 3206.4 -> 2888.6, I61: 3508 -> 3110.75, K61: 3508 -> 3110.75, M61: -3296 -> -2898.75
Set-LeveRow $wb "ARM" 61 43999 @{ H = 2888.6; I = 3110.75; K = 3110.75; M = -2898.75 }

# ARM!74 (Leve Item ID 44000): H74: 1992.4286 -> 1727.3, I74: 1991.1666 -> 1659.125, K74: 1991.1666 -> 1659.125, M74: -1117.1666 -> -785.125
Set-LeveRow $wb "ARM" 74 44000 @{ H = 1727.3; I = 1659.125; K = 1659.125; M = -785.125 }

# ARM!77 (Leve Item ID 44000): H77: 1992.4286 -> 1727.3, I77: 1991.1666 -> 1659.125, K77: 9955.833000000001 -> 8295.625, M77: -5587.833000000001 -> -3927.625
Set-LeveRow $wb "ARM" 77 44000 @{ H = 1727.3; I = 1659.125; K = 8295.625; M = -3927.625 }

# ARM!122 (Leve Item ID 36168): H122: 2949 -> 0, I122: 2949 -> 0, K122: 8847 -> 0, M122: -6397 -> ABSENT
Set-LeveRow $wb "ARM" 122 36168 @{ H = 0; I = 0; K = 0; M = $null }

# ARM!136 (Leve Item ID 43999): H136: 3206.4 -> 2888.6, I136: 3508 -> 3110.75, K136: 10524 -> 9332.25, M136: -7974 -> -6782.25
Set-LeveRow $wb "ARM" 136 43999 @{ H = 2888.6; I = 3110.75; K = 9332.25; M = -6782.25 }

# BSM!8 (Leve Item ID 2507): H8: 1263 -> 640.1667, J8: 4002.5 -> 3005, L8: 4002.5 -> 3005, N8: -4282.5 -> -3285
Set-LeveRow $wb "BSM" 8 2507 @{ H = 640.1667; J = 3005; L = 3005; N = -3285 }

# BSM!10 (Leve Item ID 2417): H10: 312.85715 -> 428, I10: 750 -> 630, J10: 138 -> 226, K10: 750 -> 630, L10: 138 -> 226, M10: -610 -> -490, N10: -418 -> -506
Set-LeveRow $wb "BSM" 10 2417 @{ H = 428; I = 630; J = 226; K = 630; L = 226; M = -490; N = -506 }

# BSM!63 (Leve Item ID 10592): H63: 30000 -> 0, J63: 30000 -> 0, L63: 30000 -> 0, N63: -31372 -> ABSENT
Set-LeveRow $wb "BSM" 63 10592 @{ H = 0; J = 0; L = 0; N = $null }

# BSM!66 (Leve Item ID 10592): H66: 30000 -> 0, J66: 30000 -> 0, L66: 90000 -> 0, N66: -96864 -> ABSENT
Set-LeveRow $wb "BSM" 66 10592 @{ H = 0; J = 0; L = 0; N = $null }

# BSM!134 (Leve Item ID 43998): H134: 1050 -> 1044, I134: 900 -> 888, K134: 2700 -> 2664, M134: -165 -> -129
Set-LeveRow $wb "BSM" 134 43998 @{ H = 1044; I = 888; K = 2664; M = -129 }

# CRP!2 (Leve Item ID 1820): H2: 1057.5555 -> 1162.1818, I2: 639.125 -> 597.6667, J2: 4405 -> 3702.5, K2: 639.125 -> 597.6667, L2: 4405 -> 3702.5, M2: -526.125 -> -484.6667, N2: -4631 -> -3928.5
Set-LeveRow $wb "CRP" 2 1820 @{ H = 1162.1818; I = 597.6667; J = 3702.5; K = 597.6667; L = 3702.5; M = -484.6667; N = -3928.5 }

# CRP!3 (Leve Item ID 3763): H3: 0 -> 1000, I3: 0 -> 1000, K3: 0 -> 1000, M3: ABSENT -> -887
Set-LeveRow $wb "CRP" 3 3763 @{ H = 1000; I = 1000; K = 1000; M = -887 }

# CRP!13 (Leve Item ID 1996): H13: 0 -> 10000, J13: 0 -> 10000, L13: 0 -> 10000, N13: ABSENT -> -10278
Set-LeveRow $wb "CRP" 13 1996 @{ H = 10000; J = 10000; L = 10000; N = -10278 }

# CRP!86 (Leve Item ID 12584): H86: 1000000000 -> 333335680, I86: 1000000000 -> 500000500, J86: 0 -> 6000, K86: 1000000000 -> 500000500, L86: 0 -> 6000, M86: -999998877 -> -499999377, N86: ABSENT -> -8246
Set-LeveRow $wb "CRP" 86 12584 @{ H = 333335680; I = 500000500; J = 6000; K = 500000500; L = 6000; M = -499999377; N = -8246 }

# CRP!89 (Leve Item ID 12584): H89: 1000000000 -> 333335680, I89: 1000000000 -> 500000500, J89: 0 -> 6000, K89: 5000000000 -> 2500002500, L89: 0 -> 30000, M89: -4999994384 -> -2499996884, N89: ABSENT -> -41232
Set-LeveRow $wb "CRP" 89 12584 @{ H = 333335680; I = 500000500; J = 6000; K = 2500002500; L = 30000; M = -2499996884; N = -41232 }

# CRP!99 (Leve Item ID 36198): H99: 1001848.4 -> 910943.25, I99: 716156.7 -> 626873.6, K99: 716156.7 -> 626873.6, M99: -714658.7 -> -625375.6
Set-LeveRow $wb "CRP" 99 36198 @{ H = 910943.25; I = 626873.6; K = 626873.6; M = -625375.6 }

# CRP!107 (Leve Item ID 27689): H107: 630.6 -> 585.2222, I107: 595.3333 -> 539.875, K107: 595.3333 -> 539.875, M107: 1324.6667 -> 1380.125
Set-LeveRow $wb "CRP" 107 27689 @{ H = 585.2222; I = 539.875; K = 539.875; M = 1380.125 }

# CRP!126 (Leve Item ID 36198): H126: 1001848.4 -> 910943.25, I126: 716156.7 -> 626873.6, K126: 2148470.1 -> 1880620.8, M126: -2146000.1 -> -1878150.8
Set-LeveRow $wb "CRP" 126 36198 @{ H = 910943.25; I = 626873.6; K = 1880620.8; M = -1878150.8 }

# CRP!132 (Leve Item ID 44019): H132: 0 -> 1700, J132: 0 -> 1700, L132: 0 -> 5100, N132: ABSENT -> -10160
Set-LeveRow $wb "CRP" 132 44019 @{ H = 1700; J = 1700; L = 5100; N = -10160 }

# CRP!134 (Leve Item ID 44020): H134: 2535.75 -> 2535.625, I134: 2547.6667 -> 2547.5, K134: 7643.000100000001 -> 7642.5, M134: -5108.000100000001 -> -5107.5
Set-LeveRow $wb "CRP" 134 44020 @{ H = 2535.625; I = 2547.5; K = 7642.5; M = -5107.5 }

# CUL!26 (Leve Item ID 4746): H26: 762 -> 809.4, J26: 2000 -> 1499.5, L26: 6000 -> 4498.5, N26: -6576 -> -5074.5
Set-LeveRow $wb "CUL" 26 4746 @{ H = 809.4; J = 1499.5; L = 4498.5; N = -5074.5 }

# CUL!131 (Leve Item ID 36060): H131: 1638.2667 -> 1901.9231, I131: 767.9 -> 978.75, K131: 2303.7 -> 2936.25, M131: 2736.3 -> 2103.75
Set-LeveRow $wb "CUL" 131 36060 @{ H = 1901.9231; I = 978.75; K = 2936.25; M = 2103.75 }

# GSM!4 (Leve Item ID 2056): H4: 0 -> 2999, I4: 0 -> 2999, K4: 0 -> 2999, M4: ABSENT -> -2887
Set-LeveRow $wb "GSM" 4 2056 @{ H = 2999; I = 2999; K = 2999; M = -2887 }

# GSM!17 (Leve Item ID 2445): H17: 1474 -> 504.33334, I17: 0 -> 6.5, J17: 1474 -> 1500, K17: 0 -> 6.5, L17: 1474 -> 1500, M17: ABSENT -> 161.5, N17: -1810 -> -1836
Set-LeveRow $wb "GSM" 17 2445 @{ H = 504.33334; I = 6.5; J = 1500; K = 6.5; L = 1500; M = 161.5; N = -1836 }

# GSM!18 (Leve Item ID 4309): H18: 279.375 -> 337.5, I18: 279.375 -> 337.5, K18: 279.375 -> 337.5, M18: 13.625 -> -44.5
Set-LeveRow $wb "GSM" 18 4309 @{ H = 337.5; I = 337.5; K = 337.5; M = -44.5 }

# GSM!122 (Leve Item ID 36182): H122: 9000 -> 8666.333000000001, J122: 8000 -> 7999.5, L122: 24000 -> 23998.5, N122: -28900 -> -28898.5
Set-LeveRow $wb "GSM" 122 36182 @{ H = 8666.333000000001; J = 7999.5; L = 23998.5; N = -28898.5 }

# GSM!126 (Leve Item ID 36184): H126: 2005.5 -> 1999, I126: 2005.5 -> 1999, K126: 6016.5 -> 5997, M126: -3546.5 -> -3527
Set-LeveRow $wb "GSM" 126 36184 @{ H = 1999; I = 1999; K = 5997; M = -3527 }

# GSM!138 (Leve Item ID 42325): H138: 0 -> 150000, J138: 0 -> 150000, L138: 0 -> 150000, N138: ABSENT -> -160280
Set-LeveRow $wb "GSM" 138 42325 @{ H = 150000; J = 150000; L = 150000; N = -160280 }

# LTW!13 (Leve Item ID 3546): H13: 7500700 -> 0, I13: 7500700 -> 0, K13: 7500700 -> 0, M13: -7500560 -> ABSENT
Set-LeveRow $wb "LTW" 13 3546 @{ H = 0; I = 0; K = 0; M = $null }

# LTW!30 (Leve Item ID 1688): H30: 828.75 -> 824.5, I30: 828.75 -> 824.5, K30: 828.75 -> 824.5, M30: -720.75 -> -716.5
Set-LeveRow $wb "LTW" 30 1688 @{ H = 824.5; I = 824.5; K = 824.5; M = -716.5 }

# LTW!46 (Leve Item ID 5282): H46: 997 -> 610.5, I46: 997 -> 610.5, K46: 997 -> 610.5, M46: -809 -> -422.5
Set-LeveRow $wb "LTW" 46 5282 @{ H = 610.5; I = 610.5; K = 610.5; M = -422.5 }

# LTW!122 (Leve Item ID 36247): H122: 4862.9473 -> 4237.3125, I122: 3316.5 -> 3350, J122: 5576.6924 -> 4533.0835, K122: 9949.5 -> 10050, L122: 16730.0772 -> 13599.2505, M122: -7499.5 -> -7600, N122: -21630.0772 -> -18499.2505
Set-LeveRow $wb "LTW" 122 36247 @{ H = 4237.3125; I = 3350; J = 4533.0835; K = 10050; L = 13599.2505; M = -7600; N = -18499.2505 }

# LTW!132 (Leve Item ID 44058): H132: 599 -> 2275, I132: 599 -> 450, J132: 0 -> 3005, K132: 1797 -> 1350, L132: 0 -> 9015, M132: 733 -> 1180, N132: ABSENT -> -14075
Set-LeveRow $wb "LTW" 132 44058 @{ H = 2275; I = 450; J = 3005; K = 1350; L = 9015; M = 1180; N = -14075 }

# LTW!136 (Leve Item ID 44060): H136: 3279 -> 3079, I136: 4000 -> 3500, K136: 12000 -> 10500, M136: -9450 -> -7950
Set-LeveRow $wb "LTW" 136 44060 @{ H = 3079; I = 3500; K = 10500; M = -7950 }

# WVR!3 (Leve Item ID 3309): H3: 3500 -> 0, I3: 3500 -> 0, K3: 3500 -> 0, M3: -3386 -> ABSENT
Set-LeveRow $wb "WVR" 3 3309 @{ H = 0; I = 0; K = 0; M = $null }

# WVR!92 (Leve Item ID 18088): H92: 0 -> 29999, J92: 0 -> 29999, L92: 0 -> 29999, N92: ABSENT -> -34991
Set-LeveRow $wb "WVR" 92 18088 @{ H = 29999; J = 29999; L = 29999; N = -34991 }
